# Commit the separate df_pontiac sheet edits:
# - Remove embedded newlines from a few "Security" column values
# - Fill in missing "Information not available" text for several Vibe rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "GTA only"
$ws.Range("E16").Value = "PL1 (GT only*)"

$ws.Range("E27").Value = "Information not available"
$ws.Range("E31").Value = "Information not available"
$ws.Range("E35").Value = "Information not available"
$ws.Range("E40").Value = "Information not available"
$ws.Range("E45").Value = "Information not available"

$ws.Range("E61").Value = "Toyota Immobilizer"
$ws.Range("E70").Value = "Toyota Immobilizer"
